$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force text to preserve formatting like "69.650.95" or "1.00"
# without Excel auto-converting to numeric values.
$dUpdates = @{
    'D2' = '69.650.95'
    'D3' = '3.805.63'
    'D4' = '1.00'
    'D5' = '613.89'
    'D6' = '177.17'
    'D7' = '3.804.51'
    'D13' = '39.74'
    'D15' = '4.435.21'
    'D16' = '3.801.79'
    'D17' = '69.710.28'
    'D21' = '506.78'
    'D22' = '9.61'
    'D25' = '86.30'
    'D32' = '8.03'
    'D33' = '31.49'
    'D35' = '0.999'
    'D39' = '481.25'
    'D40' = '0.338'
    'D41' = '3.01'
    'D43' = '49.76'
    'D44' = '44.07'
    'D45' = '8.56'
    'D46' = '2.929.25'
    'D47' = '0.0361'
    'D48' = '139.58'
    'D50' = '27.12'
    'D51' = '2.44'
}

foreach ($cellRef in $dUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $dUpdates[$cellRef]
    $cell.Style = $origStyle
}

# Column E (Volume 1h) updates - plain text with padding spaces, safe to assign directly.
$eUpdates = @{
    'E2' = '  -0.61%  '
    'E3' = '  +0.93%  '
    'E4' = '  +0.00%  '
    'E5' = '  -1.67%  '
    'E6' = '  -1.74%  '
    'E7' = '  +0.97%  '
    'E8' = '  +0.00%  '
    'E9' = '  -1.07%  '
    'E10' = '  -1.51%  '
    'E11' = '  +2.31%  '
    'E12' = '  -1.09%  '
    'E14' = '  -2.53%  '
    'E15' = '  +1.09%  '
    'E16' = '  +0.92%  '
    'E17' = '  -0.60%  '
    'E18' = '  -1.00%  '
    'E20' = '  -0.91%  '
    'E21' = '  -0.18%  '
    'E22' = '  +1.13%  '
    'E23' = '  +0.81%  '
    'E24' = '  -2.20%  '
    'E25' = '  -1.05%  '
    'E26' = '  +3.88%  '
    'E27' = '  -4.37%  '
    'E28' = '  -6.20%  '
    'E29' = '  -0.15%  '
    'E30' = '  +0.77%  '
    'E31' = '  +0.76%  '
    'E32' = '  +1.57%  '
    'E33' = '  +0.15%  '
    'E34' = '  -1.82%  '
    'E35' = '  -0.04%  '
    'E36' = '  -2.17%  '
    'E37' = '  -1.84%  '
    'E38' = '  +5.64%  '
    'E39' = '  +13.16%  '
    'E40' = '  +1.16%  '
    'E41' = '  +5.64%  '
    'E42' = '  -2.86%  '
    'E43' = '  -1.67%  '
    'E44' = '  -2.38%  '
    'E45' = '  -2.01%  '
    'E46' = '  -2.71%  '
    'E47' = '  -1.55%  '
    'E48' = '  +0.86%  '
    'E49' = '  +0.04%  '
    'E50' = '  -0.98%  '
    'E51' = '  -3.31%  '
}

foreach ($cellRef in $eUpdates.Keys) {
    $ws.Range($cellRef).Value = $eUpdates[$cellRef]
}
